$d = $word.ActiveDocument

# The last paragraph in the document currently reads "Dia 16/09: 2hr (1 dia)".
# Add a new paragraph right after it, re-using the same paragraph/run
# formatting (Arial 12pt, 360 auto line spacing, justified), with the
# text "Dia 17/09: 1hr (1 dia)".
$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.Collapse(0)  # wdCollapseEnd
$r.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "Dia 17/09: 1hr (1 dia)"
